# Mejoras en los mensajes y listado de ingreso
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the file name
$ws.Name = "dni_tutores"

# Remove the custom column width on column A
$ws.Columns.Item(1).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# Add the missing final row (duplicate of the last DNI entry)
$ws.Range("A456").Value = 41550112

# Apply left alignment to all the DNI cells; this merges/replaces the
# prior (empty) alignment + number-format styles into a single style.
$ws.Range("A1:A456").HorizontalAlignment = -4131

# Update the view: scroll back up and select the whole used range
$ws.Range("A13").Select()
$ws.Range("A1:A456").Select()
